$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column E header + data. Cells are written in this specific order so
# that new shared-string entries get created in the same order as the target
# workbook (Excel appends new unique strings to the shared string table in
# first-use order).
$ws.Range("E1").Value = "Beveilingseisen "
$ws.Range("E15").Value = "Tegen SQL en XSS beveiligd"
$ws.Range("E16").Value = "geen"
$ws.Range("E14").Value = "Elke user kan maar 1x rapporteren"
$ws.Range("E10").Value = "Comment request limit"
$ws.Range("E3").Value = "Geen"
$ws.Range("E4").Value = "Email validation"
$ws.Range("E2").Value = "Zoek voor valid ssid"
$ws.Range("E5").Value = "Geen"
$ws.Range("E6").Value = "Geen"
$ws.Range("E7").Value = "Geen"
$ws.Range("E8").Value = "Geen"
$ws.Range("E9").Value = "Geen"
$ws.Range("E11").Value = "Geen"
$ws.Range("E12").Value = "Geen"
$ws.Range("E13").Value = "Geen"
$ws.Range("E17").Value = "Tegen SQL en XSS beveiligd"

# Resize columns D and E to match the new layout.
$ws.Columns.Item(4).ColumnWidth = 17.8
$ws.Columns.Item(5).ColumnWidth = 25.7

# Update the active selection to match the saved view state.
$ws.Range("F4").Select()
